$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "1:0"
$ws.Range("H8").Value = "1:1"
$ws.Range("H9").Value = "0:0"
$ws.Range("H10").Value = "1:0"
$ws.Range("H11").Value = "1:0"
